$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 99.45
$ws.Range("K2").Value = 11.13
$ws.Range("N2").Value = 8.65

$ws.Range("K3").Value = 75
$ws.Range("N3").Value = 50

$ws.Range("G4").Value = 99.70999999999999
$ws.Range("K4").Value = 8.609999999999999
$ws.Range("N4").Value = 5.99

$ws.Range("G5").Value = 83.33
$ws.Range("K5").Value = 25
$ws.Range("N5").Value = 16.67
